$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.907.49"
$ws.Range("E2").Value = "  +3.40%  "
$ws.Range("D3").Value = "1.883.85"
$ws.Range("E3").Value = "  +3.33%  "
$ws.Range("D5").Value = "326.20"
$ws.Range("E5").Value = "  -0.58%  "
$ws.Range("D6").Value = "1.004"
$ws.Range("E6").Value = "  +0.26%  "
$ws.Range("D7").Value = "0.4683"
$ws.Range("E7").Value = "  +1.21%  "
$ws.Range("D8").Value = "0.3950"
$ws.Range("E8").Value = "  +2.63%  "
$ws.Range("D9").Value = "0.07945"
$ws.Range("E9").Value = "  +1.26%  "
$ws.Range("D10").Value = "0.9804"
$ws.Range("E10").Value = "  +2.13%  "
$ws.Range("D11").Value = "22.43"
$ws.Range("E11").Value = "  +2.86%  "
$ws.Range("D12").Value = "1.909.82"
$ws.Range("E12").Value = "  +6.98%  "
$ws.Range("D13").Value = "5.762"
$ws.Range("E13").Value = "  +2.23%  "
$ws.Range("D14").Value = "7.028"
$ws.Range("E14").Value = "  +2.58%  "
$ws.Range("D15").Value = "0.06986"
$ws.Range("E15").Value = "  +1.89%  "
$ws.Range("D16").Value = "88.73"
$ws.Range("E16").Value = "  +2.40%  "
$ws.Range("D17").Value = "1.005"
$ws.Range("E17").Value = "  +0.31%  "
$ws.Range("D18").Value = "0.00001012"
$ws.Range("E18").Value = "  +1.82%  "
$ws.Range("D19").Value = "17.02"
$ws.Range("E19").Value = "  +2.09%  "
$ws.Range("D20").Value = "1.004"
$ws.Range("E20").Value = "  +0.28%  "
$ws.Range("D21").Value = "28.893.79"
$ws.Range("E21").Value = "  +3.34%  "
$ws.Range("D22").Value = "5.374"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D23").Value = "11.15"
$ws.Range("E23").Value = "  +1.92%  "
$ws.Range("D24").Value = "2.122"
$ws.Range("E24").Value = "  +0.94%  "
$ws.Range("D25").Value = "2.104.90"
$ws.Range("E25").Value = "  +2.27%  "
$ws.Range("D26").Value = "153.63"
$ws.Range("E26").Value = "  +0.89%  "
$ws.Range("D27").Value = "19.46"
$ws.Range("E27").Value = "  +1.40%  "
$ws.Range("D28").Value = "5.788"
$ws.Range("E28").Value = "  +1.48%  "
$ws.Range("D29").Value = "2.012"
$ws.Range("E29").Value = "  +2.28%  "
$ws.Range("D30").Value = "120.22"
$ws.Range("E30").Value = "  +3.08%  "
$ws.Range("D31").Value = "0.09416"
$ws.Range("E31").Value = "  +1.76%  "
$ws.Range("D32").Value = "0.9457"
$ws.Range("E32").Value = "  +0.94%  "
$ws.Range("D33").Value = "5.328"
$ws.Range("E33").Value = "  +1.15%  "
$ws.Range("D34").Value = "1.360"
$ws.Range("E34").Value = "  +3.63%  "
$ws.Range("D35").Value = "3.352"
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D36").Value = "0.05929"
$ws.Range("E36").Value = "  -0.11%  "
$ws.Range("D37").Value = "0.02126"
$ws.Range("E37").Value = "  -1.02%  "
$ws.Range("D38").Value = "1.150"
$ws.Range("E38").Value = "  +0.24%  "
$ws.Range("D39").Value = "7.936"
$ws.Range("E39").Value = "  +4.81%  "
$ws.Range("D40").Value = "0.5728"
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("D41").Value = "10.03"
$ws.Range("E41").Value = "  +1.33%  "
$ws.Range("D42").Value = "0.1799"
$ws.Range("E42").Value = "  +1.77%  "
$ws.Range("D43").Value = "0.07269"
$ws.Range("E43").Value = "  +3.83%  "
$ws.Range("D45").Value = "0.5360"
$ws.Range("E45").Value = "  +2.21%  "
$ws.Range("D47").Value = "2.129"
$ws.Range("E47").Value = "  -3.80%  "
$ws.Range("D48").Value = "1.857"
$ws.Range("E48").Value = "  +1.87%  "
$ws.Range("D49").Value = "114.54"
$ws.Range("E49").Value = "  +1.66%  "
$ws.Range("D50").Value = "2.371"
$ws.Range("E50").Value = "  +3.14%  "
$ws.Range("E44").Value = "  +2.82%  "
$ws.Range("E46").Value = "  -7.38%  "
$ws.Range("E51").Value = "  +0.37%  "
